$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = -0.582594537033641
$ws.Cells.Item(2, 3).Value = -0.4869194962242832
$ws.Cells.Item(2, 4).Value = -0.747338533223006
$ws.Cells.Item(2, 5).Value = 0.0881164629354852
$ws.Cells.Item(2, 6).Value = -0.01196289696713632
$ws.Cells.Item(2, 7).Value = -0.3657494030468326
$ws.Cells.Item(2, 8).Value = -0.2175720126143872
$ws.Cells.Item(2, 9).Value = -0.06012602361023223
$ws.Cells.Item(2, 10).Value = -0.560785480136303
$ws.Cells.Item(2, 11).Value = -0.2804275996008339

$ws.Cells.Item(3, 2).Value = 0.09567504080935779
$ws.Cells.Item(3, 3).Value = -0.164743996189365
$ws.Cells.Item(3, 4).Value = 0.6707109999691262
$ws.Cells.Item(3, 5).Value = 0.5706316400665047
$ws.Cells.Item(3, 6).Value = 0.2168451339868084
$ws.Cells.Item(3, 7).Value = 0.3650225244192538
$ws.Cells.Item(3, 8).Value = 0.5224685134234088
$ws.Cells.Item(3, 9).Value = 0.02180905689733798
$ws.Cells.Item(3, 10).Value = 0.3021669374328071
$ws.Cells.Item(3, 11).Value = 0.1301920795877955

$ws.Cells.Item(4, 2).Value = -0.2604190369987228
$ws.Cells.Item(4, 3).Value = 0.5750359591597685
$ws.Cells.Item(4, 4).Value = 0.4749565992571469
$ws.Cells.Item(4, 5).Value = 0.1211700931774507
$ws.Cells.Item(4, 6).Value = 0.269347483609896
$ws.Cells.Item(4, 7).Value = 0.426793472614051
$ws.Cells.Item(4, 8).Value = -0.07386598391201982
$ws.Cells.Item(4, 9).Value = 0.2064918966234494
$ws.Cells.Item(4, 10).Value = 0.0345170387784377
$ws.Cells.Item(4, 11).Value = 0.3403795785247692

$ws.Cells.Item(5, 2).Value = 0.8354549961584912
$ws.Cells.Item(5, 3).Value = 0.7353756362558697
$ws.Cells.Item(5, 4).Value = 0.3815891301761735
$ws.Cells.Item(5, 5).Value = 0.5297665206086188
$ws.Cells.Item(5, 6).Value = 0.6872125096127738
$ws.Cells.Item(5, 7).Value = 0.186553053086703
$ws.Cells.Item(5, 8).Value = 0.4669109336221722
$ws.Cells.Item(5, 9).Value = 0.2949360757771605
$ws.Cells.Item(5, 10).Value = 0.600798615523492
$ws.Cells.Item(5, 11).Value = -0.01153133716379418

$ws.Cells.Item(6, 2).Value = -0.1000793599026215
$ws.Cells.Item(6, 3).Value = -0.4538658659823178
$ws.Cells.Item(6, 4).Value = -0.3056884755498724
$ws.Cells.Item(6, 5).Value = -0.1482424865457174
$ws.Cells.Item(6, 6).Value = -0.6489019430717882
$ws.Cells.Item(6, 7).Value = -0.3685440625363191
$ws.Cells.Item(6, 8).Value = -0.5405189203813308
$ws.Cells.Item(6, 9).Value = -0.2346563806349992
$ws.Cells.Item(6, 10).Value = -0.8469863333222853
$ws.Cells.Item(6, 11).Value = -0.1586149481231739

$ws.Cells.Item(7, 2).Value = -0.3537865060796963
$ws.Cells.Item(7, 3).Value = -0.2056091156472509
$ws.Cells.Item(7, 4).Value = -0.04816312664309591
$ws.Cells.Item(7, 5).Value = -0.5488225831691667
$ws.Cells.Item(7, 6).Value = -0.2684647026336975
$ws.Cells.Item(7, 7).Value = -0.4404395604787092
$ws.Cells.Item(7, 8).Value = -0.1345770207323777
$ws.Cells.Item(7, 9).Value = -0.7469069734196638
$ws.Cells.Item(7, 10).Value = -0.05853558822055238
$ws.Cells.Item(7, 11).Value = -0.3352267436446591

$ws.Cells.Item(8, 2).Value = 0.1481773904324453
$ws.Cells.Item(8, 3).Value = 0.3056233794366003
$ws.Cells.Item(8, 4).Value = -0.1950360770894705
$ws.Cells.Item(8, 5).Value = 0.08532180344599868
$ws.Cells.Item(8, 6).Value = -0.08665305439901295
$ws.Cells.Item(8, 7).Value = 0.2192094853473185
$ws.Cells.Item(8, 8).Value = -0.3931204673399676
$ws.Cells.Item(8, 9).Value = 0.2952509178591439
$ws.Cells.Item(8, 10).Value = 0.01855976243503714

$ws.Cells.Item(9, 2).Value = 0.157445989004155
$ws.Cells.Item(9, 3).Value = -0.3432134675219158
$ws.Cells.Item(9, 4).Value = -0.06285558698644665
$ws.Cells.Item(9, 5).Value = -0.2348304448314583
$ws.Cells.Item(9, 6).Value = 0.0710320949148732
$ws.Cells.Item(9, 7).Value = -0.541297857772413
$ws.Cells.Item(9, 8).Value = 0.1470735274266985
$ws.Cells.Item(9, 9).Value = -0.1296176279974082

$ws.Cells.Item(10, 2).Value = -0.5006594565260708
$ws.Cells.Item(10, 3).Value = -0.2203015759906016
$ws.Cells.Item(10, 4).Value = -0.3922764338356133
$ws.Cells.Item(10, 5).Value = -0.0864138940892818
$ws.Cells.Item(10, 6).Value = -0.698743846776568
$ws.Cells.Item(10, 7).Value = -0.01037246157745647
$ws.Cells.Item(10, 8).Value = -0.2870636170015632

$ws.Cells.Item(11, 2).Value = 0.2803578805354692
$ws.Cells.Item(11, 3).Value = 0.1083830226904575
$ws.Cells.Item(11, 4).Value = 0.414245562436789
$ws.Cells.Item(11, 5).Value = -0.1980843902504972
$ws.Cells.Item(11, 6).Value = 0.4902869949486143
$ws.Cells.Item(11, 7).Value = 0.2135958395245076

$ws.Cells.Item(12, 2).Value = -0.1719748578450117
$ws.Cells.Item(12, 3).Value = 0.1338876819013198
$ws.Cells.Item(12, 4).Value = -0.4784422707859664
$ws.Cells.Item(12, 5).Value = 0.2099291144131452
$ws.Cells.Item(12, 6).Value = -0.06676204101096155

$ws.Cells.Item(13, 2).Value = 0.3058625397463315
$ws.Cells.Item(13, 3).Value = -0.3064674129409547
$ws.Cells.Item(13, 4).Value = 0.3819039722581568
$ws.Cells.Item(13, 5).Value = 0.1052128168340501

$ws.Cells.Item(14, 2).Value = -0.6123299526872862
$ws.Cells.Item(14, 3).Value = 0.07604143251182532
$ws.Cells.Item(14, 4).Value = -0.2006497229122814

$ws.Cells.Item(15, 2).Value = 0.6883713851991116
$ws.Cells.Item(15, 3).Value = 0.4116802297750048

$ws.Cells.Item(16, 2).Value = -0.2766911554241067
